# Auto-generated edit script applying odds updates per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 1.8
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 5.5
$ws.Range("L4").Value = 5.5
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("X4").Value = 7
$ws.Range("AH4").Value = 10
$ws.Range("AI4").Value = 23
$ws.Range("AJ4").Value = 19
$ws.Range("AK4").Value = 51
$ws.Range("AM4").Value = 51
$ws.Range("AO4").Value = 10
$ws.Range("AQ4").Value = 41
$ws.Range("AX4").Value = 29
$ws.Range("BA4").Value = 151
$ws.Range("G5").Value = 2.9
$ws.Range("I5").Value = 2.35
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 1.85
$ws.Range("U5").Value = 1.75
$ws.Range("V5").Value = 2
$ws.Range("X5").Value = 15
$ws.Range("AB5").Value = 29
$ws.Range("AR5").Value = 67
$ws.Range("AS5").Value = 151
$ws.Range("S6").Value = 1.29
$ws.Range("T6").Value = 3.28
$ws.Range("G7").Value = 2.85
$ws.Range("J7").Value = 3.4
$ws.Range("L7").Value = 2.92
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 10.5
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 3.4
$ws.Range("Q7").Value = 1.65
$ws.Range("R7").Value = 1.98
$ws.Range("S7").Value = 1.33
$ws.Range("T7").Value = 3.13
$ws.Range("U7").Value = 1.5
$ws.Range("V7").Value = 2.25
$ws.Range("W7").Value = 10.5
$ws.Range("AA7").Value = 23
$ws.Range("AB7").Value = 26
$ws.Range("AC7").Value = 11.25
$ws.Range("AE7").Value = 11
$ws.Range("AH7").Value = 10
$ws.Range("AI7").Value = 14
$ws.Range("AK7").Value = 28
$ws.Range("AL7").Value = 18
$ws.Range("AM7").Value = 22
$ws.Range("AN7").Value = 5
$ws.Range("AO7").Value = 15.5
$ws.Range("AP7").Value = 20
$ws.Range("AQ7").Value = 70
$ws.Range("AT7").Value = 2.85
$ws.Range("AU7").Value = 6.1
$ws.Range("AV7").Value = 45
$ws.Range("AW7").Value = 4.5
$ws.Range("AX7").Value = 12.5
$ws.Range("AY7").Value = 17
$ws.Range("BA7").Value = 65
$ws.Range("BB7").Value = 175
$ws.Range("G11").Value = 2.05
$ws.Range("I11").Value = 3.7
$ws.Range("L11").Value = 4.5
$ws.Range("N11").Value = 7.5
$ws.Range("Q11").Value = 2.35
$ws.Range("R11").Value = 1.57
$ws.Range("S11").Value = 1.5
$ws.Range("T11").Value = 2.5
$ws.Range("AN11").Value = 4
$ws.Range("AT11").Value = 2.5
$ws.Range("K13").Value = 2.38
$ws.Range("L13").Value = 8.5
$ws.Range("Q13").Value = 1.9
$ws.Range("R13").Value = 1.95
$ws.Range("U13").Value = 2.25
$ws.Range("V13").Value = 1.57
$ws.Range("W13").Value = 6
$ws.Range("AA13").Value = 13
$ws.Range("AK13").Value = 126
$ws.Range("BA13").Value = 251
